# ggh summary table update
#
# 1. Rename the two worksheets.
# 2. Make "Growth Population & Households " the active sheet/tab (was
#    "Dwelling Types" / now "Growth by Dwelling Types").
# 3. Clear the redundant "Total Inner Ring" / "Total Outer Ring" labels
#    from column C of the subtotal rows on sheet 1.
# 4. Update the remembered selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsPop = $wb.Worksheets.Item("Population & Households by Regi")
$wsDwell = $wb.Worksheets.Item("Dwelling Types")

$wsPop.Name = "Growth Population & Households "
$wsDwell.Name = "Growth by Dwelling Types"

# Clear the duplicate "Total Inner Ring" / "Total Outer Ring" text that
# used to sit next to the subtotal rows.
$wsPop.Range("C10").ClearContents()
$wsPop.Range("C25").ClearContents()

# Update remembered selections to match the new state.
$wsPop.Range("H12").Select()
$wsDwell.Range("J6").Select()

# "Growth Population & Households " becomes the active / visible tab.
$wsPop.Activate()
